$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 26 (shifts existing rows 26.. down by one)
$ws.Rows("26:26").Insert()

# Fill in the new row 26 data: Dia=25, total_venda=20180.46, Mes=6, Ano=2025, Periodo="06/2025"
$ws.Range("A26").Value = 25
$ws.Range("B26").Value = 20180.46
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 2025
$ws.Range("E26").Value = "06/2025"

$wb.Save()
